$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) cells: force text storage so values like "26.042.71"
# or "1.00" are preserved exactly as text (matching the source inlineStr cells)
# instead of being auto-parsed as numbers by Excel, then restore the default
# "Normal" style so no stray style index/quote-prefix is left on the cell.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.042.71'
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.638.53'
$ws.Range("D3").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.0627'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '18.67'
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.752.00'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.21'
$ws.Range("D13").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.063.80'
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '191.39'
$ws.Range("D19").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '9.63'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.15'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.132'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '143.89'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.77'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.01'
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '15.26'
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0486'
$ws.Range("D30").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.878'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.128.87'
$ws.Range("D36").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0156'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '98.98'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.785'
$ws.Range("D41").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '55.50'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0527'
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.415'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.58'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.00'
$ws.Range("D49").Style = "Normal"

# Columns B (Coin), C (Link) and E (Volume(1h)) cells: plain text assignment.
$ws.Range("E2").Value = '  +0.48%  '
$ws.Range("E3").Value = '  +0.08%  '
$ws.Range("E4").Value = '  +0.56%  '
$ws.Range("E5").Value = '  -0.34%  '
$ws.Range("E6").Value = '  -0.23%  '
$ws.Range("E7").Value = '  +0.55%  '
$ws.Range("E8").Value = '  -1.60%  '
$ws.Range("E9").Value = '  -1.46%  '
$ws.Range("E10").Value = '  -4.78%  '
$ws.Range("E11").Value = '  +0.18%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("E12").Value = '  +7.13%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("E13").Value = '  -1.62%  '
$ws.Range("E14").Value = '  -2.14%  '
$ws.Range("E15").Value = '  -0.78%  '
$ws.Range("E16").Value = '  -1.93%  '
$ws.Range("E17").Value = '  +0.55%  '
$ws.Range("E18").Value = '  +0.51%  '
$ws.Range("E19").Value = '  -0.76%  '
$ws.Range("E20").Value = '  -1.58%  '
$ws.Range("E21").Value = '  -2.79%  '
$ws.Range("E22").Value = '  -1.70%  '
$ws.Range("E23").Value = '  +1.53%  '
$ws.Range("E24").Value = '  +0.30%  '
$ws.Range("B25").Value = 'Toncoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("E25").Value = '  -0.70%  '
$ws.Range("B26").Value = 'BinanceUSD'
$ws.Range("C26").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("E26").Value = '  +0.61%  '
$ws.Range("E27").Value = '  -1.35%  '
$ws.Range("E28").Value = '  -1.86%  '
$ws.Range("E29").Value = '  -0.08%  '
$ws.Range("E30").Value = '  -2.93%  '
$ws.Range("E31").Value = '  -2.29%  '
$ws.Range("E32").Value = '  -3.22%  '
$ws.Range("E33").Value = '  -1.13%  '
$ws.Range("E34").Value = '  -0.79%  '
$ws.Range("E35").Value = '  -2.47%  '
$ws.Range("E36").Value = '  -0.26%  '
$ws.Range("E37").Value = '  -0.31%  '
$ws.Range("E38").Value = '  -2.90%  '
$ws.Range("E39").Value = '  -0.79%  '
$ws.Range("E41").Value = '  -1.50%  '
$ws.Range("E42").Value = '  -3.04%  '
$ws.Range("E43").Value = '  -0.79%  '
$ws.Range("E44").Value = '  -1.98%  '
$ws.Range("E45").Value = '  -0.30%  '
$ws.Range("E46").Value = '  +1.18%  '
$ws.Range("E47").Value = '  +0.09%  '
$ws.Range("E48").Value = '  -1.16%  '
$ws.Range("E49").Value = '  +0.12%  '
$ws.Range("E50").Value = '  -2.91%  '
$ws.Range("E51").Value = '  -0.31%  '
